$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Edit the existing note in C3 (append a precision)
$ws.Range("C3").Value = "connexion a la BDD faites dans l'index (possiblement a deplacer)"

# New note in C7
$ws.Range("C7").Value = "que linterface donner ?"

# New note in C4
$ws.Range("C4").Value = "automatisation de la connexion"

# New note in C5
$ws.Range("C5").Value = "mis le projet sur GIT"

# New note in C6 (new row)
$ws.Range("C6").Value = "gant complété"
